# Apply the Sun Mar  3 21:19:08 UTC 2024 cryptos-list refresh (GitHub Actions bot)
# to the Price (D) and Volume(1h) (E) columns, plus the Stellar/TheGraph row swap
# (rows 42-43) in columns B/C/D/E.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-detected as a number by Excel
# (e.g. "415.03", "1.00") are temporarily forced to Text format so the literal
# string is preserved (matching the original inline-string cell contents), then
# restored to the default "Normal" style so no stray formatting is introduced.


# Row 2
$ws.Range('D2').Value = '62.699.67'
$ws.Range('E2').Value = '  +1.13%  '

# Row 3
$ws.Range('D3').Value = '3.466.48'
$ws.Range('E3').Value = '  +1.56%  '

# Row 4
$ws.Range('E4').Value = '  +0.07%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '415.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.37%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '129.75'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.22%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.628'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.83%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('D8').Style = 'Normal'

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.726'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.09%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.153'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +8.10%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '42.55'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.62%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.79'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +5.65%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000226'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.64%  '

# Row 14
$ws.Range('D14').Value = '4.021.97'
$ws.Range('E14').Value = '  +1.78%  '

# Row 15
$ws.Range('E15').Value = '  -0.29%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '20.56'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.22%  '

# Row 17
$ws.Range('D17').Value = '3.462.28'
$ws.Range('E17').Value = '  +0.65%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.63'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.84%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.07'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.20%  '

# Row 20
$ws.Range('D20').Value = '62.664.50'
$ws.Range('E20').Value = '  +1.14%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '462.11'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.82%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '90.32'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.89%  '

# Row 23
$ws.Range('E23').Value = '  +2.06%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.26'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.04%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.77'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +14.60%  '

# Row 26
$ws.Range('E26').Value = '  +0.25%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '33.35'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.41%  '

# Row 28
$ws.Range('E28').Value = '  +0.03%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.59'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.12%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '12.17'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.10%  '

# Row 31
$ws.Range('E31').Value = '  -1.13%  '

# Row 32
$ws.Range('E32').Value = '  -1.41%  '

# Row 33
$ws.Range('E33').Value = '  -1.91%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '40.97'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.46%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.01%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '58.30'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +8.30%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0489'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.78%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.23%  '

# Row 39
$ws.Range('E39').Value = '  +3.85%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '149.36'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.63%  '

# Row 41
$ws.Range('E41').Value = '  +5.81%  '

# Row 42
$ws.Range('B42').Value = 'TheGraph'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.323'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.21%  '

# Row 43
$ws.Range('B43').Value = 'Stellar'
$ws.Range('C43').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.135'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.33%  '

# Row 44
$ws.Range('E44').Value = '  -1.29%  '

# Row 45
$ws.Range('E45').Value = '  +4.31%  '

# Row 46
$ws.Range('E46').Value = '  +2.75%  '

# Row 47
$ws.Range('D47').Value = '0.0₃0571'
$ws.Range('E47').Value = '  +35.22%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.38'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +11.37%  '

# Row 49
$ws.Range('E49').Value = '  -1.61%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.19'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -1.39%  '

# Row 51
$ws.Range('E51').Value = '  -4.96%  '
